$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A46").Value = 44075
$ws.Range("B46").Value = 3
$ws.Range("C46").Value = "Landing page css: footer ja content"

$ws.Range("F59").Select()
